$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 111469947
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("K10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("Q10").Value = 554660.8096201464
$ws.Range("R10").Value = 6698001.275046931

# Row 11
$ws.Range("A11").Value = 111469963
$ws.Range("B11").Value = 5113
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 100526
$ws.Range("F11").Value = "Bronshjon"
$ws.Range("G11").Value = "Callidium coriaceum"
$ws.Range("H11").Value = "Paykull, 1800"
$ws.Range("K11").Value = ""
$ws.Range("M11").Value = "färska gnagspår"
$ws.Range("Q11").Value = 554718.6790950731
$ws.Range("R11").Value = 6698003.135367867

# Row 12
$ws.Range("A12").Value = 111469964
$ws.Range("B12").Value = 5113
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 100526
$ws.Range("F12").Value = "Bronshjon"
$ws.Range("G12").Value = "Callidium coriaceum"
$ws.Range("H12").Value = "Paykull, 1800"
$ws.Range("K12").Value = ""
$ws.Range("M12").Value = "färska gnagspår"
$ws.Range("Q12").Value = 554731.9372321201
$ws.Range("R12").Value = 6698141.169601779

# Row 13
$ws.Range("A13").Value = 111469949
$ws.Range("B13").Value = 96348
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = "Knärot"
$ws.Range("G13").Value = "Goodyera repens"
$ws.Range("H13").Value = "(L.) R. Br."
$ws.Range("K13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("Q13").Value = 554654.1362404823
$ws.Range("R13").Value = 6697984.37715952

# Row 14
$ws.Range("A14").Value = 111469965
$ws.Range("B14").Value = 5113
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 100526
$ws.Range("F14").Value = "Bronshjon"
$ws.Range("G14").Value = "Callidium coriaceum"
$ws.Range("H14").Value = "Paykull, 1800"
$ws.Range("K14").Value = ""
$ws.Range("M14").Value = "färska gnagspår"
$ws.Range("Q14").Value = 554716.1509068209
$ws.Range("R14").Value = 6698137.967376946

# Row 15
$ws.Range("A15").Value = 111469922
$ws.Range("B15").Value = 5135
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 105930
$ws.Range("F15").Value = "Vågbandad barkbock"
$ws.Range("G15").Value = "Semanotus undatus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("K15").Value = ""
$ws.Range("M15").Value = "äldre gnagspår"
$ws.Range("Q15").Value = 554716.6256586342
$ws.Range("R15").Value = 6698008.044787553

# Row 16
$ws.Range("A16").Value = 111469953
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("K16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("Q16").Value = 554668.8331894471
$ws.Range("R16").Value = 6698027.085862564

# Row 17
$ws.Range("A17").Value = 111469954
$ws.Range("B17").Value = 96348
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."
$ws.Range("K17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("Q17").Value = 554709.4759112563
$ws.Range("R17").Value = 6698022.75809369

# Row 18
$ws.Range("A18").Value = 111469958
$ws.Range("B18").Value = 89621
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 1101
$ws.Range("F18").Value = "Gropticka"
$ws.Range("G18").Value = "Postia guttulata"
$ws.Range("H18").Value = "(Peck) Jülich"
$ws.Range("K18").Value = ""
$ws.Range("M18").Value = ""
$ws.Range("Q18").Value = 554681.1975678616
$ws.Range("R18").Value = 6698060.372405332

# Row 19
$ws.Range("A19").Value = 111469950
$ws.Range("B19").Value = 96348
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = "Knärot"
$ws.Range("G19").Value = "Goodyera repens"
$ws.Range("H19").Value = "(L.) R. Br."
$ws.Range("K19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("Q19").Value = 554648.2514272946
$ws.Range("R19").Value = 6697980.830233379

# Row 21
$ws.Range("A21").Value = 111469941
$ws.Range("B21").Value = 96348
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 220787
$ws.Range("F21").Value = "Knärot"
$ws.Range("G21").Value = "Goodyera repens"
$ws.Range("H21").Value = "(L.) R. Br."
$ws.Range("K21").Value = "blomning"
$ws.Range("M21").Value = ""
$ws.Range("Q21").Value = 554704.8063610581
$ws.Range("R21").Value = 6698102.720679003

# Row 22
$ws.Range("A22").Value = 111469962
$ws.Range("B22").Value = 5113
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 100526
$ws.Range("F22").Value = "Bronshjon"
$ws.Range("G22").Value = "Callidium coriaceum"
$ws.Range("H22").Value = "Paykull, 1800"
$ws.Range("K22").Value = ""
$ws.Range("M22").Value = "färska gnagspår"
$ws.Range("Q22").Value = 554640.2091243146
$ws.Range("R22").Value = 6697989.107814683

# Row 23
$ws.Range("A23").Value = 111469952
$ws.Range("B23").Value = 96348
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 220787
$ws.Range("F23").Value = "Knärot"
$ws.Range("G23").Value = "Goodyera repens"
$ws.Range("H23").Value = "(L.) R. Br."
$ws.Range("K23").Value = ""
$ws.Range("M23").Value = ""
$ws.Range("Q23").Value = 554701.1291447466
$ws.Range("R23").Value = 6697985.57934437

# Row 24
$ws.Range("A24").Value = 111469926
$ws.Range("B24").Value = 89369
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 5447
$ws.Range("F24").Value = "Vedticka"
$ws.Range("G24").Value = "Fuscoporia viticola"
$ws.Range("H24").Value = "(Schwein.) Murrill"
$ws.Range("K24").Value = ""
$ws.Range("M24").Value = ""
$ws.Range("Q24").Value = 554745.7538377594
$ws.Range("R24").Value = 6698078.142900295

# Row 25
$ws.Range("A25").Value = 111469951
$ws.Range("B25").Value = 96348
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("K25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("Q25").Value = 554679.0891228422
$ws.Range("R25").Value = 6697970.425878088

# Row 26
$ws.Range("A26").Value = 111469968
$ws.Range("B26").Value = 5113
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 100526
$ws.Range("F26").Value = "Bronshjon"
$ws.Range("G26").Value = "Callidium coriaceum"
$ws.Range("H26").Value = "Paykull, 1800"
$ws.Range("K26").Value = ""
$ws.Range("M26").Value = "äldre gnagspår"
$ws.Range("Q26").Value = 554679.218646974
$ws.Range("R26").Value = 6698060.342582431

# Row 28
$ws.Range("A28").Value = 111469966
$ws.Range("B28").Value = 5113
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 100526
$ws.Range("F28").Value = "Bronshjon"
$ws.Range("G28").Value = "Callidium coriaceum"
$ws.Range("H28").Value = "Paykull, 1800"
$ws.Range("K28").Value = ""
$ws.Range("M28").Value = "äldre gnagspår"
$ws.Range("Q28").Value = 554729.2459973614
$ws.Range("R28").Value = 6698057.144588907

# Row 29
$ws.Range("A29").Value = 111469967
$ws.Range("B29").Value = 5113
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 100526
$ws.Range("F29").Value = "Bronshjon"
$ws.Range("G29").Value = "Callidium coriaceum"
$ws.Range("H29").Value = "Paykull, 1800"
$ws.Range("K29").Value = ""
$ws.Range("M29").Value = "äldre gnagspår"
$ws.Range("Q29").Value = 554705.6319759471
$ws.Range("R29").Value = 6698113.601669285

# Row 30
$ws.Range("A30").Value = 111469946
$ws.Range("B30").Value = 96348
$ws.Range("D30").Value = "VU"
$ws.Range("E30").Value = 220787
$ws.Range("F30").Value = "Knärot"
$ws.Range("G30").Value = "Goodyera repens"
$ws.Range("H30").Value = "(L.) R. Br."
$ws.Range("K30").Value = ""
$ws.Range("M30").Value = ""
$ws.Range("Q30").Value = 554664.6782300239
$ws.Range("R30").Value = 6698007.261790544

# Row 31
$ws.Range("A31").Value = 111469969
$ws.Range("B31").Value = 76495
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 6487
$ws.Range("F31").Value = "Blågrå svartspik"
$ws.Range("G31").Value = "Chaenothecopsis fennica"
$ws.Range("H31").Value = "(Laurila) Tibell"
$ws.Range("K31").Value = ""
$ws.Range("M31").Value = ""
$ws.Range("Q31").Value = 554769.2275642991
$ws.Range("R31").Value = 6698129.381786803

# Row 32
$ws.Range("A32").Value = 111469940
$ws.Range("B32").Value = 79444
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 1049
$ws.Range("F32").Value = "Kortskaftad ärgspik"
$ws.Range("G32").Value = "Microcalicium ahlneri"
$ws.Range("H32").Value = "Tibell"
$ws.Range("K32").Value = ""
$ws.Range("M32").Value = ""
$ws.Range("Q32").Value = 554771.1915359092
$ws.Range("R32").Value = 6698130.399477887
